# Add "Penthouse Rental" to the itinerary (column E) and budget (column F)
# on row 5, matching the style used by the existing F4 "$"#,##0.00 entries
# in the row below it, and update the selected cell to reflect where the
# user last clicked after making the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New itinerary item name (adds a new shared string: "Penthouse Rental")
$ws.Range("E5").Value = "Penthouse Rental"

# New budget amount for the penthouse rental
$ws.Range("F5").Value = 8323
$ws.Range("F5").NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'

# Update the active selection, mirroring the author's final cursor position
$ws.Range("F16").Select()
